# Updated cryptos list - applies Price (D) and Volume(1h) (E) changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.010.62'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").Value = '2.412.61'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.52'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.22'
$ws.Range("E6").Value = '  -0.70%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +0.65%  '

$ws.Range("E9").Value = '  -1.21%  '

$ws.Range("E10").Value = '  -2.28%  '

$ws.Range("E11").Value = '  -0.57%  '

$ws.Range("E12").Value = '  -1.26%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.73'
$ws.Range("E13").Value = '  -0.14%  '

$ws.Range("D14").Value = '2.840.92'
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").Value = '59.897.76'
$ws.Range("E15").Value = '  -0.01%  '

$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("D17").Value = '2.403.89'
$ws.Range("E17").Value = '  -0.49%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.18'
$ws.Range("E18").Value = '  -1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  +3.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '327.82'
$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.74'
$ws.Range("E21").Value = '  +0.69%  '

$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.77'
$ws.Range("E23").Value = '  -1.66%  '

$ws.Range("E24").Value = '  +3.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.61'
$ws.Range("E25").Value = '  +0.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.40'
$ws.Range("E27").Value = '  +4.29%  '

$ws.Range("E28").Value = '  +1.19%  '

$ws.Range("E29").Value = '  -1.48%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.84'
$ws.Range("E30").Value = '  +0.09%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.13'
$ws.Range("E31").Value = '  -0.98%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.12'
$ws.Range("E32").Value = '  +9.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.401'
$ws.Range("E33").Value = '  -2.87%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.43'
$ws.Range("E34").Value = '  -0.94%  '

$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("E36").Value = '  +2.66%  '

$ws.Range("E37").Value = '  -0.02%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '322.16'
$ws.Range("E39").Value = '  +2.55%  '

$ws.Range("E40").Value = '  -0.67%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.76'
$ws.Range("E41").Value = '  +6.09%  '

$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.81'
$ws.Range("E44").Value = '  +2.67%  '

$ws.Range("E45").Value = '  -0.99%  '

$ws.Range("E46").Value = '  -0.56%  '

$ws.Range("E47").Value = '  -1.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '11.04'
$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("E49").Value = '  -1.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.66'
$ws.Range("E50").Value = '  -0.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.938'
$ws.Range("E51").Value = '  -2.14%  '
